$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.477.50"
$ws.Range("E2").Value = "  -1.12%  "

$ws.Range("D3").Value = "1.618.28"
$ws.Range("E3").Value = "  -1.94%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "211.27"
$ws.Range("E5").Value = "  -1.09%  "

$ws.Range("D6").Value = "0.525"
$ws.Range("E6").Value = "  -1.50%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "22.83"
$ws.Range("E8").Value = "  -1.64%  "

$ws.Range("E9").Value = "  +0.43%  "

$ws.Range("E10").Value = "  -0.56%  "

$ws.Range("E11").Value = "  -0.46%  "

$ws.Range("D12").Value = "1.847.27"
$ws.Range("E12").Value = "  -1.86%  "

$ws.Range("D13").Value = "1.619.70"
$ws.Range("E13").Value = "  -1.74%  "

$ws.Range("E14").Value = "  -0.44%  "

$ws.Range("E15").Value = "  -2.68%  "

$ws.Range("D16").Value = "65.01"
$ws.Range("E16").Value = "  +0.71%  "

$ws.Range("D17").Value = "27.463.75"
$ws.Range("E17").Value = "  -1.03%  "

$ws.Range("D18").Value = "233.12"
$ws.Range("E18").Value = "  +0.06%  "

$ws.Range("D20").Value = "7.54"
$ws.Range("E20").Value = "  -2.38%  "

$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").Value = "4.30"
$ws.Range("E22").Value = "  -0.73%  "

$ws.Range("D23").Value = "10.19"
$ws.Range("E23").Value = "  +0.21%  "

$ws.Range("E24").Value = "  +5.51%  "

$ws.Range("D25").Value = "150.53"
$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("E26").Value = "  -1.77%  "

$ws.Range("E27").Value = "  -1.46%  "

$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("D29").Value = "15.54"
$ws.Range("E29").Value = "  -1.05%  "

$ws.Range("E30").Value = "  -1.24%  "

$ws.Range("E31").Value = "  -1.18%  "

$ws.Range("E32").Value = "  -1.39%  "

$ws.Range("D33").Value = "1.471.14"
$ws.Range("E33").Value = "  +2.07%  "

$ws.Range("E34").Value = "  -3.21%  "

$ws.Range("E35").Value = "  -3.28%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "2.33"
$ws.Range("E36").Value = "  -0.49%  "

$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "0.953"
$ws.Range("E37").Value = "  +8.25%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.0167"
$ws.Range("E38").Value = "  -0.52%  "

$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "0.558"
$ws.Range("E39").Value = "  -2.68%  "

$ws.Range("E40").Value = "  -2.99%  "

$ws.Range("E41").Value = "  +0.07%  "

$ws.Range("D42").Value = "67.99"
$ws.Range("E42").Value = "  +2.04%  "

$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("D44").Value = "0.980"
$ws.Range("E44").Value = "  -5.15%  "

$ws.Range("E45").Value = "  -2.50%  "

$ws.Range("D46").Value = "5.25"
$ws.Range("E46").Value = "  -6.66%  "

$ws.Range("D47").Value = "1.758.20"
$ws.Range("E47").Value = "  -1.86%  "

$ws.Range("D48").Value = "1.73"
$ws.Range("E48").Value = "  +0.48%  "

$ws.Range("D49").Value = "86.82"
$ws.Range("E49").Value = "  +0.27%  "

$ws.Range("E50").Value = "  -2.41%  "

$ws.Range("E51").Value = "  +1.19%  "
